$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1358.2307
$ws.Range("I19").Value = 1482.2858
$ws.Range("J19").Value = 1213.5
$ws.Range("K19").Value = 1482.2858
$ws.Range("L19").Value = 1213.5
$ws.Range("M19").Value = -1307.2858
$ws.Range("N19").Value = -1563.5
$ws.Range("H32").Value = 1999
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1999
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1999
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2651
$ws.Range("H33").Value = 210.75
$ws.Range("I33").Value = 180.55173
$ws.Range("K33").Value = 180.55173
$ws.Range("M33").Value = 48.44827000000001
$ws.Range("H40").Value = 6642.7144
$ws.Range("I40").Value = 11999.833
$ws.Range("J40").Value = 2624.875
$ws.Range("K40").Value = 11999.833
$ws.Range("L40").Value = 2624.875
$ws.Range("M40").Value = -11824.833
$ws.Range("N40").Value = -2974.875
$ws.Range("H51").Value = 10000
$ws.Range("J51").Value = 10000
$ws.Range("L51").Value = 10000
$ws.Range("N51").Value = -10968
$ws.Range("H69").Value = 142864860
$ws.Range("I69").Value = 333338000
$ws.Range("J69").Value = 10000
$ws.Range("K69").Value = 1000014000
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -1000013126
$ws.Range("N69").Value = -31748
$ws.Range("H72").Value = 142864860
$ws.Range("I72").Value = 333338000
$ws.Range("J72").Value = 10000
$ws.Range("K72").Value = 3000042000
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -3000037632
$ws.Range("N72").Value = -98736
$ws.Range("H74").Value = 8993
$ws.Range("I74").Value = 8993
$ws.Range("K74").Value = 8993
$ws.Range("M74").Value = -8057
$ws.Range("H77").Value = 8993
$ws.Range("I77").Value = 8993
$ws.Range("K77").Value = 44965
$ws.Range("M77").Value = -40285
$ws.Range("H86").Value = 4464.5
$ws.Range("I86").Value = 4448
$ws.Range("J86").Value = 4472.75
$ws.Range("K86").Value = 4448
$ws.Range("L86").Value = 4472.75
$ws.Range("M86").Value = -3325
$ws.Range("N86").Value = -6718.75
$ws.Range("H88").Value = 1652.6
$ws.Range("J88").Value = 1680.6666
$ws.Range("L88").Value = 1680.6666
$ws.Range("N88").Value = -2492.6666
$ws.Range("H89").Value = 4464.5
$ws.Range("I89").Value = 4448
$ws.Range("J89").Value = 4472.75
$ws.Range("K89").Value = 22240
$ws.Range("L89").Value = 22363.75
$ws.Range("M89").Value = -16624
$ws.Range("N89").Value = -33595.75
$ws.Range("H91").Value = 1652.6
$ws.Range("J91").Value = 1680.6666
$ws.Range("L91").Value = 1680.6666
$ws.Range("N91").Value = -4488.6666
$ws.Range("H125").Value = 4141.5713
$ws.Range("J125").Value = 4598.8
$ws.Range("L125").Value = 41389.2
$ws.Range("N125").Value = -46309.2
$ws.Range("H127").Value = 940
$ws.Range("I127").Value = 940
$ws.Range("K127").Value = 2820
$ws.Range("M127").Value = 2140
$ws.Range("H131").Value = 696.1429000000001
$ws.Range("J131").Value = 115
$ws.Range("L131").Value = 345
$ws.Range("N131").Value = -10425
$ws.Range("H132").Value = 557622.9399999999
$ws.Range("I132").Value = 1814.2
$ws.Range("K132").Value = 5442.6
$ws.Range("M132").Value = -2912.6
$ws.Range("H138").Value = 3779.3
$ws.Range("I138").Value = 2512
$ws.Range("J138").Value = 4888.1875
$ws.Range("K138").Value = 7536
$ws.Range("L138").Value = 14664.5625
$ws.Range("M138").Value = -2396
$ws.Range("N138").Value = -24944.5625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2468.4375
$ws.Range("I2").Value = 2399.7144
$ws.Range("K2").Value = 2399.7144
$ws.Range("M2").Value = -2286.7144
$ws.Range("H5").Value = 2666.3333
$ws.Range("J5").Value = 2999.5
$ws.Range("L5").Value = 2999.5
$ws.Range("N5").Value = -3223.5
$ws.Range("H32").Value = 22180.162
$ws.Range("I32").Value = 4275.28
$ws.Range("K32").Value = 4275.28
$ws.Range("M32").Value = -3988.28
$ws.Range("H34").Value = 25732
$ws.Range("J34").Value = 26172.924
$ws.Range("L34").Value = 26172.924
$ws.Range("N34").Value = -26714.924
$ws.Range("H45").Value = 1163.4166
$ws.Range("I45").Value = 1163.4166
$ws.Range("K45").Value = 1163.4166
$ws.Range("M45").Value = -786.4166
$ws.Range("H61").Value = 2811.0715
$ws.Range("I61").Value = 2334.2856
$ws.Range("K61").Value = 2334.2856
$ws.Range("M61").Value = -2122.2856
$ws.Range("H74").Value = 1201.6
$ws.Range("I74").Value = 1212.2106
$ws.Range("K74").Value = 1212.2106
$ws.Range("M74").Value = -338.2106000000001
$ws.Range("H77").Value = 1201.6
$ws.Range("I77").Value = 1212.2106
$ws.Range("K77").Value = 6061.053000000001
$ws.Range("M77").Value = -1693.053000000001
$ws.Range("H116").Value = 2468.4375
$ws.Range("I116").Value = 2399.7144
$ws.Range("K116").Value = 2399.7144
$ws.Range("M116").Value = -105.7143999999998
$ws.Range("H132").Value = 2178.75
$ws.Range("I132").Value = 825.7826
$ws.Range("K132").Value = 2477.3478
$ws.Range("M132").Value = 52.65219999999999
$ws.Range("H136").Value = 2811.0715
$ws.Range("I136").Value = 2334.2856
$ws.Range("K136").Value = 7002.8568
$ws.Range("M136").Value = -4452.8568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2468.4375
$ws.Range("I3").Value = 2399.7144
$ws.Range("K3").Value = 2399.7144
$ws.Range("M3").Value = -2285.7144
$ws.Range("H4").Value = 2666.3333
$ws.Range("J4").Value = 2999.5
$ws.Range("L4").Value = 2999.5
$ws.Range("N4").Value = -3229.5
$ws.Range("H12").Value = 3933
$ws.Range("I12").Value = 10000
$ws.Range("K12").Value = 10000
$ws.Range("M12").Value = -9832
$ws.Range("H20").Value = 1940.35
$ws.Range("I20").Value = 2231
$ws.Range("K20").Value = 2231
$ws.Range("M20").Value = -1984
$ws.Range("H22").Value = 501.2143
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 103564.836
$ws.Range("J58").Value = 103564.836
$ws.Range("L58").Value = 103564.836
$ws.Range("N58").Value = -104152.836
$ws.Range("H81").Value = 44779.4
$ws.Range("J81").Value = 45724.5
$ws.Range("L81").Value = 45724.5
$ws.Range("N81").Value = -47846.5
$ws.Range("H82").Value = 19249.2
$ws.Range("I82").Value = 3123.5
$ws.Range("K82").Value = 3123.5
$ws.Range("M82").Value = -2740.5
$ws.Range("H84").Value = 44779.4
$ws.Range("J84").Value = 45724.5
$ws.Range("L84").Value = 137173.5
$ws.Range("N84").Value = -147781.5
$ws.Range("H85").Value = 19249.2
$ws.Range("I85").Value = 3123.5
$ws.Range("K85").Value = 3123.5
$ws.Range("M85").Value = -1797.5
$ws.Range("H86").Value = 1900
$ws.Range("I86").Value = 1900
$ws.Range("K86").Value = 1900
$ws.Range("M86").Value = -777
$ws.Range("H89").Value = 1900
$ws.Range("I89").Value = 1900
$ws.Range("K89").Value = 9500
$ws.Range("M89").Value = -3884
$ws.Range("H99").Value = 1952.3636
$ws.Range("J99").Value = 2999.75
$ws.Range("L99").Value = 2999.75
$ws.Range("N99").Value = -5995.75
$ws.Range("H104").Value = 89992.5
$ws.Range("J104").Value = 89992.5
$ws.Range("L104").Value = 89992.5
$ws.Range("N104").Value = -96980.5
$ws.Range("H105").Value = 4462.357
$ws.Range("I105").Value = 4665.3335
$ws.Range("J105").Value = 4407
$ws.Range("K105").Value = 4665.3335
$ws.Range("L105").Value = 4407
$ws.Range("M105").Value = -2918.3335
$ws.Range("N105").Value = -7901
$ws.Range("H107").Value = 1574.8334
$ws.Range("I107").Value = 1462.6666
$ws.Range("K107").Value = 1462.6666
$ws.Range("M107").Value = 457.3334
$ws.Range("H134").Value = 2511.5652
$ws.Range("I134").Value = 1962.1765
$ws.Range("J134").Value = 4068.1667
$ws.Range("K134").Value = 5886.529500000001
$ws.Range("L134").Value = 12204.5001
$ws.Range("M134").Value = -3351.529500000001
$ws.Range("N134").Value = -17274.5001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1500
$ws.Range("I14").Value = 1500
$ws.Range("K14").Value = 1500
$ws.Range("M14").Value = -1330
$ws.Range("H31").Value = 4718.695
$ws.Range("I31").Value = 3162
$ws.Range("K31").Value = 3162
$ws.Range("M31").Value = -2867
$ws.Range("H34").Value = 4718.695
$ws.Range("I34").Value = 3162
$ws.Range("K34").Value = 3162
$ws.Range("M34").Value = -2960
$ws.Range("H39").Value = 19925.334
$ws.Range("I39").Value = 8888
$ws.Range("K39").Value = 8888
$ws.Range("M39").Value = -8497
$ws.Range("H49").Value = 19925.334
$ws.Range("I49").Value = 8888
$ws.Range("K49").Value = 8888
$ws.Range("M49").Value = -8706
$ws.Range("H52").Value = 88633
$ws.Range("J52").Value = 88633
$ws.Range("L52").Value = 88633
$ws.Range("N52").Value = -89221
$ws.Range("H58").Value = 5226.4443
$ws.Range("I58").Value = 7333
$ws.Range("K58").Value = 7333
$ws.Range("M58").Value = -7130
$ws.Range("H99").Value = 2722.8386
$ws.Range("I99").Value = 2885.087
$ws.Range("K99").Value = 2885.087
$ws.Range("M99").Value = -1387.087
$ws.Range("H107").Value = 566
$ws.Range("I107").Value = 547.9
$ws.Range("K107").Value = 547.9
$ws.Range("M107").Value = 1372.1
$ws.Range("H122").Value = 1933.1
$ws.Range("I122").Value = 2036.0555
$ws.Range("J122").Value = 1006.5
$ws.Range("K122").Value = 6108.166499999999
$ws.Range("L122").Value = 3019.5
$ws.Range("M122").Value = -3658.166499999999
$ws.Range("N122").Value = -7919.5
$ws.Range("H126").Value = 2722.8386
$ws.Range("I126").Value = 2885.087
$ws.Range("K126").Value = 8655.261
$ws.Range("M126").Value = -6185.261
$ws.Range("H136").Value = 5226.4443
$ws.Range("I136").Value = 7333
$ws.Range("K136").Value = 21999
$ws.Range("M136").Value = -19449

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 181.63637
$ws.Range("I14").Value = 181.63637
$ws.Range("K14").Value = 544.9091100000001
$ws.Range("M14").Value = -371.9091100000001
$ws.Range("H45").Value = 5392
$ws.Range("J45").Value = 5000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -16064
$ws.Range("H131").Value = 75762.87
$ws.Range("J131").Value = 2650.7144
$ws.Range("L131").Value = 7952.1432
$ws.Range("N131").Value = -18032.1432
$ws.Range("H136").Value = 9030.5
$ws.Range("I136").Value = 9030.5
$ws.Range("K136").Value = 27091.5
$ws.Range("M136").Value = -21991.5
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18610350
$ws.Range("J11").Value = 8262750
$ws.Range("L11").Value = 8262750
$ws.Range("N11").Value = -8263028
$ws.Range("H80").Value = 8859.9
$ws.Range("I80").Value = 5392.923
$ws.Range("J80").Value = 11511.117
$ws.Range("K80").Value = 5392.923
$ws.Range("L80").Value = 11511.117
$ws.Range("M80").Value = -4394.923
$ws.Range("N80").Value = -13507.117
$ws.Range("H83").Value = 8859.9
$ws.Range("I83").Value = 5392.923
$ws.Range("J83").Value = 11511.117
$ws.Range("K83").Value = 26964.615
$ws.Range("L83").Value = 57555.585
$ws.Range("M83").Value = -21972.615
$ws.Range("N83").Value = -67539.58499999999
$ws.Range("H97").Value = 2749.5
$ws.Range("I97").Value = 2666.1667
$ws.Range("J97").Value = 2999.5
$ws.Range("K97").Value = 2666.1667
$ws.Range("L97").Value = 2999.5
$ws.Range("M97").Value = -2170.1667
$ws.Range("N97").Value = -3991.5
$ws.Range("H102").Value = 2821.6428
$ws.Range("I102").Value = 2316.875
$ws.Range("J102").Value = 3494.6667
$ws.Range("K102").Value = 2316.875
$ws.Range("L102").Value = 3494.6667
$ws.Range("M102").Value = -694.875
$ws.Range("N102").Value = -6738.6667
$ws.Range("H132").Value = 4947.5454
$ws.Range("I132").Value = 5455.7144
$ws.Range("K132").Value = 16367.1432
$ws.Range("M132").Value = -13837.1432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10999.667
$ws.Range("I40").Value = 10999
$ws.Range("K40").Value = 10999
$ws.Range("M40").Value = -10863
$ws.Range("H55").Value = 785.8
$ws.Range("J55").Value = 732.5
$ws.Range("L55").Value = 732.5
$ws.Range("N55").Value = -1078.5
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H100").Value = 4531.3335
$ws.Range("I100").Value = 5000
$ws.Range("J100").Value = 3594
$ws.Range("K100").Value = 5000
$ws.Range("L100").Value = 3594
$ws.Range("M100").Value = -4459
$ws.Range("N100").Value = -4676
$ws.Range("H132").Value = 3300.4375
$ws.Range("J132").Value = 2917.25
$ws.Range("L132").Value = 8751.75
$ws.Range("N132").Value = -13811.75
$ws.Range("H136").Value = 3665.3845
$ws.Range("J136").Value = 4174
$ws.Range("L136").Value = 12522
$ws.Range("N136").Value = -17622
$ws.Range("H138").Value = 71250
$ws.Range("J138").Value = 71250
$ws.Range("L138").Value = 71250
$ws.Range("N138").Value = -81530

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 14497.8
$ws.Range("I8").Value = 16000
$ws.Range("J8").Value = 14122.25
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 14122.25
$ws.Range("M8").Value = -15860
$ws.Range("N8").Value = -14402.25
$ws.Range("H15").Value = 14990
$ws.Range("J15").Value = 14990
$ws.Range("L15").Value = 14990
$ws.Range("N15").Value = -15566
$ws.Range("H40").Value = 49999
$ws.Range("I40").Value = 49999
$ws.Range("K40").Value = 49999
$ws.Range("M40").Value = -49850
$ws.Range("H51").Value = 35000
$ws.Range("I51").Value = 45000
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 45000
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -44490
$ws.Range("N51").Value = -26020
$ws.Range("H52").Value = 40047
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H54").Value = 48500
$ws.Range("J54").Value = 45000
$ws.Range("L54").Value = 45000
$ws.Range("N54").Value = -46040
$ws.Range("H113").Value = 1123
$ws.Range("I113").Value = 917.8570999999999
$ws.Range("K113").Value = 2753.5713
$ws.Range("M113").Value = -583.5712999999996
$ws.Range("H126").Value = 2866.5
$ws.Range("I126").Value = 3239.8
$ws.Range("K126").Value = 9719.400000000001
$ws.Range("M126").Value = -7249.400000000001
$ws.Range("H136").Value = 3540.9092
$ws.Range("I136").Value = 3834
$ws.Range("K136").Value = 11502
$ws.Range("M136").Value = -8952

